# T1176_Companies_AddFVAOpportunityOnCompanyDetailPage.xlsx
# "Companies Changes 1 - 29th June 2023"
#
# - Users sheet: replace the test user "Nicole Bicho" with "Drew Koecher"
# - Leave AddOpportunity's selection where it was, but make Users the
#   active sheet/tab, with C9 selected
# - Tidy up the stray "apply number format" flag that a few header/value
#   cells on AddOpportunity were carrying (no visible format change,
#   General stays General) - clears to the same look while dropping the
#   redundant flag

$wb = $excel.ActiveWorkbook

$usersSheet = $wb.Worksheets.Item("Users")
$usersSheet.Range("A2").Value = "Drew Koecher"

$oppSheet = $wb.Worksheets.Item("AddOpportunity")
$oppSheet.Range("J1:K1").ClearFormats()
$oppSheet.Range("J1:K1").Font.Bold = $true
$oppSheet.Range("R1:U1").ClearFormats()
$oppSheet.Range("R1:U1").Font.Bold = $true
$oppSheet.Range("J2:K2").ClearFormats()
$oppSheet.Range("R2:U2").ClearFormats()

$usersSheet.Activate()
$usersSheet.Range("C9").Select()
